$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 3 & 4 ("Golden Arrows v Swallows" / "Supersport Utd v Richards Bay")
#    were re-ordered: swap the match data held in columns F:V while leaving
#    the row index (A), competition metadata (B:D) and kickoff date (E)
#    untouched on each row.
# ---------------------------------------------------------------------------
$row3 = $ws.Range("F3:V3").Value2
$row4 = $ws.Range("F4:V4").Value2
$ws.Range("F3:V3").Value2 = $row4
$ws.Range("F4:V4").Value2 = $row3

# ---------------------------------------------------------------------------
# 2) Rows 7 & 8 ("TS Galaxy v Cape Town Spurs" / "Kaizer Chiefs v Chippa Utd.")
#    were likewise re-ordered: swap columns F:V.
# ---------------------------------------------------------------------------
$row7 = $ws.Range("F7:V7").Value2
$row8 = $ws.Range("F8:V8").Value2
$ws.Range("F7:V7").Value2 = $row8
$ws.Range("F8:V8").Value2 = $row7

# ---------------------------------------------------------------------------
# 3) Append two new match rows (55 & 56) at the bottom of the sheet, copying
#    the formatting of the last existing data row (54) so the bold/bordered
#    index column and the date/time number format carry over correctly.
# ---------------------------------------------------------------------------
$ws.Range("A54:V54").Copy()
$ws.Range("A55:V55").PasteSpecial(-4122)
$ws.Range("A54:V54").Copy()
$ws.Range("A56:V56").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A55").Value2 = 54
$ws.Range("B55").Value2 = "south-africa"
$ws.Range("C55").Value2 = "premier-league"
$ws.Range("D55").Value2 = "2023-2024"
$ws.Range("E55").Value2 = 45196.8125
$ws.Range("F55").Value2 = "Kaizer Chiefs"
$ws.Range("G55").Value2 = 2
$ws.Range("H55").Value2 = "Sekhukhune"
$ws.Range("I55").Value2 = 1
$ws.Range("J55").Value2 = 1.85
$ws.Range("K55").Value2 = "20/09/2023 18:43"
$ws.Range("L55").Value2 = 2.15
$ws.Range("M55").Value2 = "27/09/2023 19:29"
$ws.Range("N55").Value2 = 3.15
$ws.Range("O55").Value2 = "20/09/2023 18:43"
$ws.Range("P55").Value2 = 2.89
$ws.Range("Q55").Value2 = "27/09/2023 19:29"
$ws.Range("R55").Value2 = 5.02
$ws.Range("S55").Value2 = "20/09/2023 18:43"
$ws.Range("T55").Value2 = 4.29
$ws.Range("U55").Value2 = "27/09/2023 19:29"
$ws.Range("V55").Value2 = "https://www.betexplorer.com/football/south-africa/premier-league/kaizer-chiefs-sekhukhune/vBmHAig5/"

$ws.Range("A56").Value2 = 55
$ws.Range("B56").Value2 = "south-africa"
$ws.Range("C56").Value2 = "premier-league"
$ws.Range("D56").Value2 = "2023-2024"
$ws.Range("E56").Value2 = 45196.8125
$ws.Range("F56").Value2 = "Mamelodi Sundowns"
$ws.Range("G56").Value2 = 3
$ws.Range("H56").Value2 = "Stellenbosch"
$ws.Range("I56").Value2 = 1
$ws.Range("J56").Value2 = 1.33
$ws.Range("K56").Value2 = "25/09/2023 13:12"
$ws.Range("L56").Value2 = 1.34
$ws.Range("M56").Value2 = "27/09/2023 19:21"
$ws.Range("N56").Value2 = 4.69
$ws.Range("O56").Value2 = "25/09/2023 13:12"
$ws.Range("P56").Value2 = 4.68
$ws.Range("Q56").Value2 = "27/09/2023 19:22"
$ws.Range("R56").Value2 = 10.32
$ws.Range("S56").Value2 = "25/09/2023 13:12"
$ws.Range("T56").Value2 = 10.75
$ws.Range("U56").Value2 = "27/09/2023 19:22"
$ws.Range("V56").Value2 = "https://www.betexplorer.com/football/south-africa/premier-league/mamelodi-sundowns-stellenbosch-fc/MZxM9B8B/"
